$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20 (shifts existing row 20 "ERASTAPEX..." and all
# rows below it down by one), for the new item "EMPACOZA TRIO XR 25/5/1000 30TAB"
# which sorts alphabetically right before "ERASTAPEX".
$ws.Rows.Item(20).Insert()

# Clone the formatting (borders/fill/font/number-format) of the row that now
# sits at 21 (the former row 20) onto the freshly inserted, blank row 20 so the
# new row matches the rest of the table visually.
$ws.Range("A21:Q21").Copy()
$ws.Range("A20:Q20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row's data.
$ws.Cells.Item(20, 1).Value = 14
$ws.Cells.Item(20, 3).Value = "EMPACOZA TRIO XR 25/5/1000  30TAB"
$ws.Cells.Item(20, 8).Value = "1:0"
$ws.Cells.Item(20, 12).Value = "0"
$ws.Cells.Item(20, 14).Value = "396.00"
$ws.Cells.Item(20, 16).Value = "130.6800"
$ws.Cells.Item(20, 17).Value = "0:1"

# Update the running total (shifted down to row 54 by the insert) to include
# the new item's selling price (2465.18 + 130.68 = 2595.86).
$ws.Cells.Item(54, 16).Value = 2595.86

# Update the generated-on timestamp in the footer (shifted down to row 55).
$ws.Cells.Item(55, 1).Value = "Sunday, 3 August, 2025 12:24 PM"
